$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting the old "citation" column (C) to D
$ws.Columns("C").Insert()

# Row 1 headers
$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "month"
$ws.Range("C1").Value = "method name"
$ws.Range("D1").Value = "citation"
$ws.Range("E1").Value = "multiple_instruments"
$ws.Range("F1").Value = "multiple_outcomes"

# Row 2 - existing data row: re-home method/citation from B/C into C/D, add E
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "omnigenic Mendelian randomization"
$ws.Range("D2").Value = "wang_mendelian_nodate"
$ws.Range("E2").Value = "yes"

# Row 3 - new record
$ws.Range("A3").Value = 2021
$ws.Range("C3").Value = "moPMR-Egger"
$ws.Range("D3").Value = "liu_multi-trait_2021"
$ws.Range("E3").Value = "yes"
$ws.Range("F3").Value = "yes"

# Row 4 - new record
$ws.Range("A4").Value = 2020
$ws.Range("C4").Value = "PMR-Egger"
$ws.Range("D4").Value = "yuan_testing_2020"
$ws.Range("E4").Value = "yes"
$ws.Range("F4").Value = "no"

# Column widths (closest achievable values given the engine's pixel-quantized
# column-width model; matches the target widths of 30.33 / 21.3 / 18.24 / 16.71)
$ws.Columns("C").ColumnWidth = 29.5
$ws.Columns("D").ColumnWidth = 20.5
$ws.Columns("E").ColumnWidth = 17.33
$ws.Columns("F").ColumnWidth = 15.83

# Move the active selection, matching the post-edit cursor position
[void]$ws.Range("C5").Select()
